$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function setRGB($idx, $hex) {
  $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
  $val = $r + ($g * 256) + ($b * 65536)
  $c = $tcs.Colors($idx)
  $c.RGB = $val
}

setRGB 1 "000000"
setRGB 2 "FFFFFF"
setRGB 3 "44546A"
setRGB 4 "E7E6E6"
setRGB 5 "5B9BD5"
setRGB 6 "ED7D31"
setRGB 7 "A5A5A5"
setRGB 8 "FFC000"
setRGB 9 "4472C4"
setRGB 10 "70AD47"
setRGB 11 "0563C1"
setRGB 12 "954F72"
